$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 19 (Number of models = 9, run2 numbers) ---
$ws.Range("D19").Value = 1098.6300000000001
$ws.Range("E19").Value = 1098.6300000000001
$ws.Range("F19").Value = 982.69799999999998
$ws.Range("G19").Value = 1038.06

# --- Row 20 (Number of models = 10, run2 numbers) ---
$ws.Range("D20").Value = 1126.78
$ws.Range("E20").Value = 1126.78
$ws.Range("F20").Value = 1054.03
$ws.Range("G20").Value = 1064.8800000000001

# --- Row 23 (scaling block 2, run #1, run2 numbers) ---
$ws.Range("D23").Value = 2504.63
$ws.Range("E23").Value = 2504.63
$ws.Range("F23").Value = 2298.7800000000002
$ws.Range("G23").Value = 2294.7399999999998

# --- Row 24 ---
$ws.Range("D24").Value = 2556.46
$ws.Range("E24").Value = 2556.46
$ws.Range("F24").Value = 2355.91
$ws.Range("G24").Value = 2398.65

# --- Row 25 ---
$ws.Range("D25").Value = 2610.67
$ws.Range("E25").Value = 2610.67
$ws.Range("F25").Value = 2493.2199999999998
$ws.Range("G25").Value = 2426.2399999999998

# --- Row 26 ---
$ws.Range("D26").Value = 2497.88
$ws.Range("E26").Value = 2497.88
$ws.Range("F26").Value = 2346.5700000000002
$ws.Range("G26").Value = 2339.84

# --- Row 27 ---
$ws.Range("D27").Value = 2626.53
$ws.Range("E27").Value = 2626.53
$ws.Range("F27").Value = 2509.63
$ws.Range("G27").Value = 2229.88

# --- Row 28 ---
$ws.Range("D28").Value = 2492.1999999999998
$ws.Range("E28").Value = 2492.1999999999998
$ws.Range("F28").Value = 2374.19
$ws.Range("G28").Value = 2257.34

# --- Row 29 ---
$ws.Range("D29").Value = 2514.3000000000002
$ws.Range("E29").Value = 2514.3000000000002
$ws.Range("F29").Value = 2396.31
$ws.Range("G29").Value = 2358.08

# --- Row 30 ---
$ws.Range("D30").Value = 2586.29
$ws.Range("E30").Value = 2586.29
$ws.Range("F30").Value = 2467.23
$ws.Range("G30").Value = 2330.46

# --- Row 31 ---
$ws.Range("D31").Value = 2584.27
$ws.Range("E31").Value = 2584.27
$ws.Range("F31").Value = 2430.89
$ws.Range("G31").Value = 2451.56

# --- Move the view/selection to where the new data was entered ---
$ws.Activate()
$ws.Range("D20").Select()
